{"js": "// Word JavaScript API (Office.js) script.\n// Body of: async (context) => { ... }\n//\n// Applies three content edits described by the diff:\n//   1. Fix the title typo \"SCREEING REPORT\" -> \"SCREENING REPORT\" by\n//      splitting the run so the inserted \"N\" lands in its own run and\n//      the \"_GoBack\" bookmark (Word's \"last edit location\" marker) sits\n//      right after it - exactly mirroring the authored XML.\n//   2. Move the \"_GoBack\" bookmark: remove it from the end of the\n//      document (where it previously sat after the trailing page break)\n//      since Word only ever keeps a single \"_GoBack\" bookmark, tracking\n//      the most recent edit location.\n//   3. Update the report date \"2020-03-13\" -> \"2021-03-04\".\n\nconst doc = context.document;\nconst body = doc.body;\n\n// --- Step 1: remove the stale \"_GoBack\" bookmark from the end of the doc ---\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Step 2: fix \"SCREEING REPORT\" -> \"SCREENING REPORT\", splitting runs\n//     so \"N\" is its own run and the \"_GoBack\" bookmark follows it ---\nconst titleResults = body.search(\"SCREEING REPORT\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  const titleRange = titleResults.items[0];\n\n  const titleOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr><w:b/><w:u w:val=\"single\"/></w:rPr>\n              <w:t>SCREE</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:b/><w:u w:val=\"single\"/></w:rPr>\n              <w:t>N</w:t>\n            </w:r>\n            <w:bookmarkStart w:id=\"100\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"100\"/>\n            <w:r>\n              <w:rPr><w:b/><w:u w:val=\"single\"/></w:rPr>\n              <w:t>ING REPORT</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\n  titleRange.insertOoxml(titleOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Step 3: update the report date ---\nconst dateResults = body.search(\"2020-03-13\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2021-03-04\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is already open as $d.\n#\n# Applies three content edits described by the diff:\n#   1. Remove the stale \"_GoBack\" bookmark from the end of the document\n#      (it previously sat after the trailing page break).\n#   2. Fix the title typo \"SCREEING REPORT\" -> \"SCREENING REPORT\" by\n#      splitting the run so the inserted \"N\" lands in its own run and\n#      the \"_GoBack\" bookmark (Word's \"last edit location\" marker) sits\n#      right after it - exactly mirroring the authored XML. Word only\n#      ever keeps a single \"_GoBack\" bookmark, tracking the most recent\n#      edit location, which is why it moves here.\n#   3. Update the report date \"2020-03-13\" -> \"2021-03-04\".\n\n$d = $word.ActiveDocument\n\n# --- Step 1: remove the stale \"_GoBack\" bookmark from the end of the doc ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- Step 2: fix \"SCREEING REPORT\" -> \"SCREENING REPORT\", splitting runs\n#     so \"N\" is its own run and the \"_GoBack\" bookmark follows it ---\n$titleRange = $d.Content\n$titleFound = $titleRange.Find.Execute(\"SCREEING REPORT\")\nif ($titleFound) {\n  $titleOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' +\n    '<w:r><w:rPr><w:b/><w:u w:val=\"single\"/></w:rPr><w:t>SCREE</w:t></w:r>' +\n    '<w:r><w:rPr><w:b/><w:u w:val=\"single\"/></w:rPr><w:t>N</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"100\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"100\"/>' +\n    '<w:r><w:rPr><w:b/><w:u w:val=\"single\"/></w:rPr><w:t>ING REPORT</w:t></w:r>' +\n    '</w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  $titleRange.InsertXML($titleOoxml)\n}\n\n# --- Step 3: update the report date ---\n$dateRange = $d.Content\n$dateRange.Find.Execute(\"2020-03-13\", $false, $false, $false, $false, $false, $true, 1, $false, \"2021-03-04\", 2) | Out-Null\n"}
